$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are numeric-looking strings (e.g. "585.90", "1.00")
# must be forced to Text format first so Excel keeps them as literal strings
# (matching the original inline-string cell contents) instead of auto-converting
# them to numbers and losing formatting such as trailing zeros.
$textForceCells = @(
    "D5", "D6", "D9", "D10", "D11", "D12", "D13", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D39", "D40", "D41", "D42", "D45", "D47", "D48", "D49", "D50"
)
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.122.59"
$ws.Range("E2").Value = "  -2.72%  "
$ws.Range("D3").Value = "3.131.90"
$ws.Range("E3").Value = "  -5.59%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "585.90"
$ws.Range("E5").Value = "  -2.84%  "
$ws.Range("D6").Value = "134.77"
$ws.Range("E6").Value = "  -5.60%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "3.128.85"
$ws.Range("E8").Value = "  -5.65%  "
$ws.Range("D9").Value = "0.505"
$ws.Range("E9").Value = "  -2.86%  "
$ws.Range("D10").Value = "0.140"
$ws.Range("E10").Value = "  -6.42%  "
$ws.Range("D11").Value = "5.21"
$ws.Range("E11").Value = "  -5.02%  "
$ws.Range("D12").Value = "0.451"
$ws.Range("E12").Value = "  -4.47%  "
$ws.Range("D13").Value = "0.0000232"
$ws.Range("E13").Value = "  -6.72%  "
$ws.Range("D14").Value = "33.77"
$ws.Range("E14").Value = "  -2.93%  "
$ws.Range("D15").Value = "3.651.64"
$ws.Range("E15").Value = "  -5.27%  "
$ws.Range("E16").Value = "  -2.04%  "
$ws.Range("D17").Value = "3.133.12"
$ws.Range("E17").Value = "  -5.35%  "
$ws.Range("D18").Value = "62.329.97"
$ws.Range("E18").Value = "  -2.51%  "
$ws.Range("D19").Value = "6.51"
$ws.Range("E19").Value = "  -5.42%  "
$ws.Range("D20").Value = "450.57"
$ws.Range("E20").Value = "  -6.23%  "
$ws.Range("D21").Value = "13.79"
$ws.Range("E21").Value = "  -3.07%  "
$ws.Range("D22").Value = "0.698"
$ws.Range("E22").Value = "  -4.99%  "
$ws.Range("D23").Value = "7.53"
$ws.Range("E23").Value = "  -6.82%  "
$ws.Range("D24").Value = "13.23"
$ws.Range("E24").Value = "  -4.01%  "
$ws.Range("D25").Value = "83.04"
$ws.Range("E25").Value = "  -1.96%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  -3.76%  "
$ws.Range("D29").Value = "7.63"
$ws.Range("E29").Value = "  -6.42%  "
$ws.Range("D30").Value = "6.70"
$ws.Range("E30").Value = "  -8.52%  "
$ws.Range("D31").Value = "1.99"
$ws.Range("E31").Value = "  -8.15%  "
$ws.Range("D32").Value = "27.00"
$ws.Range("E32").Value = "  -6.40%  "
$ws.Range("D33").Value = "0.102"
$ws.Range("E33").Value = "  -4.34%  "
$ws.Range("D34").Value = "2.35"
$ws.Range("E34").Value = "  -8.11%  "
$ws.Range("E35").Value = "  -8.30%  "
$ws.Range("D36").Value = "5.75"
$ws.Range("E36").Value = "  -4.49%  "
$ws.Range("D37").Value = "50.72"
$ws.Range("E37").Value = "  -5.54%  "
$ws.Range("D38").Value = "0.0₃0695"
$ws.Range("E38").Value = "  -7.03%  "
$ws.Range("D39").Value = "0.0383"
$ws.Range("E39").Value = "  -4.54%  "
$ws.Range("D40").Value = "2.65"
$ws.Range("E40").Value = "  -4.08%  "
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").Value = "7.99"
$ws.Range("E41").Value = "  -4.54%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "391.47"
$ws.Range("E42").Value = "  -9.80%  "
$ws.Range("E43").Value = "  -3.83%  "
$ws.Range("D44").Value = "2.741.85"
$ws.Range("E44").Value = "  -10.84%  "
$ws.Range("D45").Value = "0.248"
$ws.Range("E45").Value = "  -6.77%  "
$ws.Range("D47").Value = "2.10"
$ws.Range("E47").Value = "  -5.14%  "
$ws.Range("D48").Value = "124.59"
$ws.Range("E48").Value = "  -2.96%  "
$ws.Range("D49").Value = "25.00"
$ws.Range("E49").Value = "  -5.50%  "
$ws.Range("D50").Value = "34.10"
$ws.Range("E50").Value = "  -5.60%  "
$ws.Range("E51").Value = "  -4.04%  "
